$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Golang Developer", "https://www.dice.com/job-detail/4e0cad2a-5167-4096-9148-440fb179fdce", "Richmond, Virginia", "Contract, Third Party", "$50 - $60", "Lorvenk Technologies LLC"),
    @("Golang Developer", "https://www.dice.com/job-detail/fe1336e0-6e7a-4b97-a450-08f8ec0c26ba", "Virginia Beach, Virginia", "Contract", "Depends on Experience", "Oxford Global Resources"),
    @("Junior Level Golang Developer ::: Remote ::: Long Term", "https://www.dice.com/job-detail/9855f6bb-0df0-4ed8-b580-9c6d003298f3", "Remote", "Contract", "Depends on Experience", "VRN Technologies"),
    @("Mid Level Golang Developer ::: Remote ::: Long Term", "https://www.dice.com/job-detail/225c92b3-6786-4bf0-b965-2bb542eefcd8", "Remote", "Contract", "Depends on Experience", "VRN Technologies"),
    @("Senior Golang Developer ::: Remote ::: Long Term", "https://www.dice.com/job-detail/1ab059fa-aefa-4ed1-b01d-5c4adb80ad5a", "Remote", "Contract", "Depends on Experience", "ERPA"),
    @("Integration and Middleware Developer (.NET, GoLang, Python) - Harrisburg, PA - hybrid role - 12 Months", "https://www.dice.com/job-detail/4af1581d-3f55-4119-b445-59e5bb9f0392", "Harrisburg, Pennsylvania", "Contract, Third Party", "Depends on Experience", "StarTechs Inc."),
    @("Golang Architect / Principal Backend Architect", "https://www.dice.com/job-detail/d0aee673-010b-42d8-84bf-74658e016c0f", "Atlanta, Georgia", "Third Party", "$58 - $68 per hour", "Quantum World Technologies Inc."),
    @("Golang Architect / Principal Backend Architect", "https://www.dice.com/job-detail/8cc35663-69dc-4e57-8cde-e05042f54f6f", "Atlanta, Georgia", "Third Party, Contract", "$58 - $68 per hour", "Stellite Works LLC"),
    @("Lead / Sr Golang Developer", "https://www.dice.com/job-detail/21c59a65-f637-49e9-8aa6-9ffff64d94c8", "Hybrid in Plano, Texas", "Contract", "Depends on Experience", "Lorvin Technologies")
)

$startRow = 17
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
